$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the (formerly) last paragraph of
#    the document to right after the ")" that closes the registration
#    number in the letterhead ("...002774428-W)").
# ---------------------------------------------------------------------

# Remove the existing bookmark at the end of the document, if present.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete() | Out-Null
}

# Locate the ")" that closes the registration number. NB: re-using the
# same Range object after Find.Execute is important - Word collapses /
# resizes that Range to the match in place.
$rng = $d.Content
$rng.Find.Execute("002774428-W)", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$insertPos = $rng.End

# Work around an engine quirk: creating a bookmark on a zero-length
# Range that sits exactly at a paragraph-end position snaps it back to
# the start of the document. Instead, insert a temporary marker
# character right after the ")", anchor the bookmark just before that
# marker (now a perfectly ordinary, non-paragraph-end position), and
# then delete the marker again - the collapsed bookmark stays put.
$markerRange = $d.Range($insertPos, $insertPos)
$markerRange.InsertAfter("@@MARK@@")

$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.Content.Find.Execute("@@MARK@@", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Replace "DRIED SEA CUCUMBER" with "PAVAKA" (both occurrences in the
#    table).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("DRIED SEA CUCUMBER", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "PAVAKA", 2) | Out-Null
